$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 73597824
$ws.Range("B2").Value = 93276
$ws.Range("D2").Value = "LC"
$ws.Range("E2").Value = 2170
$ws.Range("F2").Value = "Flagellkvastmossa"
$ws.Range("G2").Value = "Dicranum flagellare"
$ws.Range("H2").Value = "Hedw."
$ws.Range("L2").Value = "'"
$ws.Range("L2").ClearFormats()
$ws.Range("Q2").Value = 551313.8124669526
$ws.Range("R2").Value = 6516081.021942991
$ws.Range("S2").Value = 25

$ws.Range("A3").Value = 73597823
$ws.Range("B3").Value = 90676
$ws.Range("D3").Value = "NT"
$ws.Range("E3").Value = 5966
$ws.Range("F3").Value = "Motaggsvamp"
$ws.Range("G3").Value = "Sarcodon squamosus"
$ws.Range("H3").Value = "(Schaeff.) Quél."
$ws.Range("L3").ClearContents()
$ws.Range("Q3").Value = 551313.8124669526
$ws.Range("R3").Value = 6516081.021942991
$ws.Range("S3").Value = 25

$ws.Range("A4").Value = 73597772
$ws.Range("B4").Value = 90008
$ws.Range("D4").Value = "LC"
$ws.Range("E4").Value = 6031
$ws.Range("F4").Value = "Blomkålssvamp"
$ws.Range("G4").Value = "Sparassis crispa"
$ws.Range("H4").Value = "(Wulfen:Fr.) Fr."
$ws.Range("L4").ClearContents()
$ws.Range("Q4").Value = 551214.4045088139
$ws.Range("R4").Value = 6515978.080872892
$ws.Range("S4").Value = 25

$ws.Range("A5").Value = 73597770
$ws.Range("B5").Value = 93375
$ws.Range("D5").Value = "LC"
$ws.Range("E5").Value = 2180
$ws.Range("F5").Value = "Blåmossa"
$ws.Range("G5").Value = "Leucobryum glaucum"
$ws.Range("H5").Value = "(Hedw.) Ångstr."
$ws.Range("L5").Value = "'"
$ws.Range("L5").ClearFormats()
$ws.Range("Q5").Value = 551214.4045088139
$ws.Range("R5").Value = 6515978.080872892
$ws.Range("S5").Value = 25

$ws.Range("A6").Value = 73597795
$ws.Range("B6").Value = 77177
$ws.Range("D6").Value = "NT"
$ws.Range("E6").Value = 353
$ws.Range("F6").Value = "Dvärgbägarlav"
$ws.Range("G6").Value = "Cladonia parasitica"
$ws.Range("H6").Value = "(Hoffm.) Hoffm."
$ws.Range("L6").ClearContents()
$ws.Range("Q6").Value = 551373.0161088589
$ws.Range("R6").Value = 6515925.219612571
$ws.Range("S6").Value = 50

$ws.Range("A7").Value = 73597810
$ws.Range("B7").Value = 77259
$ws.Range("D7").Value = "NT"
$ws.Range("E7").Value = 228912
$ws.Range("F7").Value = "Mörk kolflarnlav"
$ws.Range("G7").Value = "Carbonicola myrmecina"
$ws.Range("H7").Value = "(Ach.) Bendiksby & Timdal"
$ws.Range("L7").ClearContents()
$ws.Range("Q7").Value = 551373.0161088589
$ws.Range("R7").Value = 6515925.219612571
$ws.Range("S7").Value = 50
